# MP6532_board_28-QFN BOM updates:
#  - order-detail / availability notes added for a couple of parts
#  - new LCSC part numbers for replacements received from JLC
#  - a couple of designator lists grew (extra resistors/caps added to the design)
#  - column A widened so the longer "Comment" text is readable, selection moved to D7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Comment / Designator / LCSC Part # text updates -----------------------

# Row 4: 470nF/100V cap - note that the previously specified LCSC part is out of stock
$ws.Range("A4").Value = "470nF/100V (not available on order C97926 )"
$ws.Range("D4").Value = "C16195147"

# Row 6: C 100nF caps - three more designators added (C10,C11,C12)
$ws.Range("B6").Value = "C17,C19,C18,C7,C8,C9,C10,C11,C12"

# Row 11: LED RED - annotate that the old LCSC part isn't available; the note is
# appended as a second (differently-formatted) run within the same cell, same as
# the source workbook which carries the annotation as rich text.
$ws.Range("A11").Value = "LED RED (C2295 not available on order)"
$ws.Range("A11").Characters(10, 30).Font.ColorIndex = -4105
$ws.Range("D11").Value = "C965812"

# Row 13: Res 10K - three more designators added (R25,R26,R27)
$ws.Range("B13").Value = "R8,R9,R10,R11,R12,R13,R25,R26,R27"

# Row 14: MP6532 controller - record where good parts actually came from
$ws.Range("A14").Value = "MP6532 from My parts / digikey 1589-MP6532GR-ZCT-ND "

# --- Sheet view / column layout --------------------------------------------

# Column A needs to be much wider now that the Comment text is longer
$ws.Columns("A").ColumnWidth = 64.6

# Move the active selection to D7 (matches the author's last click before saving)
$ws.Range("D7").Select() | Out-Null
